# RBA v2.5 - Atualizacao da Tela
# Replace placeholder "TRE"/"Tre"/"tre" tokens with "QWER"/"Qwer"/"Qewr"/"qwer"
# tokens across the document body and the page header.
#
# wdReplaceNone = 0, wdReplaceOne = 1, wdReplaceAll = 2 (standard Word enum)
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)

$d = $word.ActiveDocument

# --- 1. Body: "A TERE," salutation (single bold occurrence) -----------------
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- 2. Header: placeholder tokens ------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)

# "DIRETORIA DE ENSINO REGIAO TRE" -> "... QWER"
$r = $hdr.Range
$r.Find.Execute("TRE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# own paragraph "TERE - DEP." -> "QWER - DEP."
$r = $hdr.Range
$r.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# address paragraph: five "Tre" occurrences get distinct replacements, in order
$treReplacements = @("Qwer", "Qwer", "Qewr", "Qewr", "Qwer")
$r = $hdr.Range
foreach ($rep in $treReplacements) {
    $r.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, $rep, 1) | Out-Null
    $r.Collapse(0)
}

# three lowercase "tre" occurrences (CEP, Tel, Email) -> "qwer"
$r = $hdr.Range
for ($i = 0; $i -lt 3; $i++) {
    $r.Find.Execute("tre", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
    $r.Collapse(0)
}
